$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." as:
#   <w:proofErr spellStart/> <r>Version</r> <proofErr spellEnd/> <r> 1.</r>
#   <bookmarkStart _GoBack/> <bookmarkEnd/>
#
# Target reads "Version 2." as:
#   <w:proofErr spellStart/> <r>Versi</r> <r>on</r> <proofErr spellEnd/> <r> 2</r>
#   <bookmarkStart _GoBack/> <bookmarkEnd/> <r>.</r>

# Step 1: split the "Version" run into "Versi" + "on" (run break added
# mid-word), keeping the spell-check proofErr markers in their correct
# positions (spellStart before "Versi", spellEnd after "on").
$r = $d.Range(5, 7)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Step 2: bump the version number, 1 -> 2.
$r2 = $d.Range(8, 9)
$r2.Text = "2"

# Step 3: pull the trailing "." out of the " 2." run ...
$r3 = $d.Range(9, 10)
$r3.Text = ""

# Step 4: ... and re-insert it as its own run after the _GoBack bookmark.
$r4 = $d.Range(10, 10)
$r4.InsertAfter(".")
